# "Move 4.0 -> 4.2"
#
# The Django docs link on the "Model Field Types" slide points at the 4.0
# version of the docs; bump it to 4.2. The run containing the literal
# "/4.0/" segment is removed and its replacement ("/4.2/") is folded into
# the following run (which keeps going straight into
# "ref/models/fields/#field-types"), exactly mirroring how PowerPoint
# re-tokenizes runs when you retype a substring in the middle of a text box.

$p = $ppt.ActivePresentation

# "Model Field Types" is slide 13 in the deck.
$slide = $p.Slides.Item(13)
$shape = $slide.Shapes.Item("Rectangle 5")
$tr = $shape.TextFrame.TextRange

$fullText = $tr.Text
$oldSegment = "/4.0/"
$newSegment = "/4.2/"

$start = $fullText.IndexOf($oldSegment)
if ($start -ge 0) {
    # 1-based character position for PowerPoint's Characters() indexer.
    $charPos = $start + 1
    $tailLength = $tr.Length - $start

    # Replace everything from the start of "/4.0/" through the end of the
    # string in one shot, so PowerPoint merges it into a single run instead
    # of leaving a separate run behind for the edited segment.
    $tail = $tr.Characters($charPos, $tailLength)
    $tail.Text = $newSegment + $fullText.Substring($start + $oldSegment.Length)
}
